# repull data, push all data, mean calculation
# Update the "dSF" (column F) values for a set of rows to reflect the
# repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    10 = -3
    20 = -3
    23 = 2
    24 = 0
    25 = 3
    30 = 1
    32 = -4
    34 = -2
    35 = 1
    38 = 2
    40 = 0
    45 = -2
    46 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
